# Applies the diff: splits the original single run of paragraph 1 into
# three runs separated by proofErr spell-check markers around "bzw", and
# appends a brand-new second paragraph describing why histograms were
# removed (also containing proofErr spell/grammar markers).

$d = $word.ActiveDocument
$p1 = $d.Paragraphs.Item(1)

$wordMlNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

$para1 = '<w:p ' + $wordMlNs + '>' +
         '<w:r><w:t xml:space="preserve">Nans </w:t></w:r>' +
         '<w:proofErr w:type="spellStart"/>' +
         '<w:r><w:t>bzw</w:t></w:r>' +
         '<w:proofErr w:type="spellEnd"/>' +
         '<w:r><w:t xml:space="preserve"> fehlende werte werden immer mit -1 ersetzt</w:t></w:r>' +
         '</w:p>'

$para2 = '<w:p ' + $wordMlNs + '>' +
         '<w:r><w:t xml:space="preserve">Histogramme entfernt </w:t></w:r>' +
         '<w:r><w:t xml:space="preserve">#Wegen absolut niedriger Feature </w:t></w:r>' +
         '<w:proofErr w:type="spellStart"/>' +
         '<w:r><w:t>Importance</w:t></w:r>' +
         '<w:proofErr w:type="spellEnd"/>' +
         '<w:r><w:t xml:space="preserve"> nach </w:t></w:r>' +
         '<w:proofErr w:type="spellStart"/>' +
         '<w:r><w:t>RandomForestFit</w:t></w:r>' +
         '<w:proofErr w:type="spellEnd"/>' +
         '<w:r><w:t xml:space="preserve"> | </w:t></w:r>' +
         '<w:proofErr w:type="gramStart"/>' +
         '<w:r><w:t>Verbraucht</w:t></w:r>' +
         '<w:proofErr w:type="gramEnd"/>' +
         '<w:r><w:t xml:space="preserve"> sehr viel Platz | Sehr viele Features</w:t></w:r>' +
         '</w:p>'

# Replacing the whole paragraph range (which includes its trailing
# paragraph mark) with two w:p fragments turns the one paragraph into two,
# without leaving a stray empty paragraph behind.
[void]$p1.Range.InsertXML($para1 + $para2)
